$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("07-09-2021", 6831588, 0),
    @("08-09-2021", 6130159, 0),
    @("09-09-2021", 3537896, 82620),
    @("10-09-2021", 2215600, 112514),
    @("13-09-2021", 3117100, 282475),
    @("14-09-2021", 3138200, 257312)
)

$startRow = 173
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
